$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("B1").Value = "Person Name"
$ws.Range("C1").Value = "Supervisors"
$ws.Range("D1").Value = "Supervisees"
$ws.Range("E1").Value = "Source URL"
$ws.Range("F1").Value = "Timestamp"
$ws.Range("G1").Value = "Notes"

# Old headers occupied H1:I1 ("Timestamp", "Notes") - the sheet shrank to
# 7 columns (A:G), so clear the now-unused trailing columns entirely.
$ws.Range("H1:I1").Clear()

# --- Row 2: update existing assessment row ---
$ws.Range("A2").Value = "e9cf0b2c-106e-465f-b52a-8641b08367ca"
$ws.Range("B2").Value = "Lisette Espin-Noboa"
$ws.Range("C2").Value = "ss"
$ws.Range("D2").Value = "ss"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "2025-09-04T04:45:01.976Z"
$ws.Range("G2").Value = "Assessment for Lisette Espin-Noboa's genealogy"

# Clear the stale H2/I2 cells that used to hold this row's Timestamp/Notes.
$ws.Range("H2:I2").Clear()

# --- Row 3: new assessment row ---
$ws.Range("A3").Value = "d6678c3d-f2ff-4b95-b841-b1b65d99c87c"
$ws.Range("B3").Value = "Albert Einstein"
$ws.Range("C3").Value = "dsas"
$ws.Range("D3").Value = "dasd"
$ws.Range("E3").Value = "sda"
$ws.Range("F3").Value = "2025-09-04T08:16:17.922Z"
$ws.Range("G3").Value = "Assessment for Albert Einstein's genealogy"

# --- Row 4: new assessment row ---
$ws.Range("A4").Value = "27c800cc-48a0-4760-bfb2-b3e57e23b685"
$ws.Range("B4").Value = "Albert Einstein"
$ws.Range("C4").Value = "sds"
$ws.Range("D4").Value = "dds"
$ws.Range("E4").Value = "https://mail.google.com/mail/u/3/#inbox/FMfcgzQcpTMQrQMqWsrPPZBdhMvFsrXs"
$ws.Range("F4").Value = "2025-09-04T08:16:42.337Z"
$ws.Range("G4").Value = "Assessment for Albert Einstein's genealogy"
